# Updated symbol list (Price / Volume(1h) / Hora columns) for rows 2-51 of
# the crypto ranking sheet, matching the latest coinranking.com pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target literal text per cell (column letter + row number -> new text).
# These must stay TEXT cells (they were written as inline strings originally),
# so a leading "'" forces Excel to store them as text instead of auto-coercing
# "310.60" / "0.37%" / "2" into numbers, and resetting the style back to
# "Normal" afterwards drops the quote-prefix formatting Excel would otherwise
# stamp on the cell, leaving the cell's own style untouched (matches source).
$cells = @{
  "D2" = "310.60"
  "E2" = "0.37%"
  "G2" = "2"
  "D3" = "40.41"
  "E3" = "-1.59%"
  "G3" = "2"
  "D4" = "5.079"
  "E4" = "-2.67%"
  "G4" = "2"
  "D5" = "0.07544"
  "E5" = "-1.72%"
  "G5" = "2"
  "D6" = "4.320"
  "E6" = "-0.34%"
  "G6" = "2"
  "D7" = "1.716"
  "E7" = "6.21%"
  "G7" = "2"
  "D8" = "0.9303"
  "E8" = "1.37%"
  "G8" = "2"
  "D9" = "2.423"
  "E9" = "-0.91%"
  "G9" = "2"
  "D10" = "0.1238"
  "E10" = "0.07%"
  "G10" = "2"
  "D11" = "0.1803"
  "E11" = "-1.58%"
  "G11" = "2"
  "D12" = "0.09085"
  "E12" = "-0.66%"
  "G12" = "2"
  "D13" = "0.04103"
  "E13" = "-5.07%"
  "G13" = "2"
  "D14" = "0.1052"
  "E14" = "0.02%"
  "G14" = "2"
  "D15" = "0.001292"
  "E15" = "2.40%"
  "G15" = "2"
  "D16" = "0.005982"
  "E16" = "2.92%"
  "G16" = "2"
  "E17" = "0.32%"
  "G17" = "2"
  "D18" = "3.348"
  "E18" = "-0.20%"
  "G18" = "2"
  "D19" = "0.3354"
  "E19" = "0.56%"
  "G19" = "2"
  "D20" = "7.693"
  "E20" = "5.97%"
  "G20" = "2"
  "D21" = "0.1356"
  "E21" = "-3.30%"
  "G21" = "2"
  "D22" = "0.2936"
  "E22" = "0.69%"
  "G22" = "2"
  "D23" = "0.04049"
  "E23" = "-0.68%"
  "G23" = "2"
  "D24" = "0.001266"
  "E24" = "0.49%"
  "G24" = "2"
  "D25" = "0.004051"
  "E25" = "-0.79%"
  "G25" = "2"
  "D26" = "0.0001274"
  "E26" = "0.14%"
  "G26" = "2"
  "G27" = "2"
  "G28" = "2"
  "G29" = "2"
  "G30" = "2"
  "G31" = "2"
  "G32" = "2"
  "G33" = "2"
  "G34" = "2"
  "G35" = "2"
  "G36" = "2"
  "G37" = "2"
  "D38" = "0.02424"
  "E38" = "-0.75%"
  "G38" = "2"
  "D39" = "0.05151"
  "E39" = "-2.39%"
  "G39" = "2"
  "D40" = "0.007726"
  "E40" = "-1.45%"
  "G40" = "2"
  "D41" = "0.1294"
  "E41" = "-1.49%"
  "G41" = "2"
  "D42" = "0.007702"
  "G42" = "2"
  "D43" = "0.002187"
  "E43" = "14.29%"
  "G43" = "2"
  "D44" = "0.008024"
  "E44" = "-3.84%"
  "G44" = "2"
  "D45" = "0.3096"
  "E45" = "-7.28%"
  "G45" = "2"
  "D46" = "0.00006645"
  "E46" = "-3.01%"
  "G46" = "2"
  "D47" = "0.00000000753"
  "E47" = "0.19%"
  "G47" = "2"
  "D48" = "0.2617"
  "E48" = "27.35%"
  "G48" = "2"
  "D49" = "0.004216"
  "E49" = "2.84%"
  "G49" = "2"
  "D50" = "0.00002108"
  "E50" = "0.19%"
  "G50" = "2"
  "D51" = "0.0002008"
  "E51" = "0.19%"
  "G51" = "2"
}

foreach ($ref in $cells.Keys) {
  $rng = $ws.Range($ref)
  $rng.Value = "'" + $cells[$ref]
  $rng.Style = "Normal"
}
